# Added New Mac-Address and Document Types
# Appends 5 new rows (157-161) to the master-reg_center_machine_device_h
# data sheet, switches the workbook to manual calculation, and updates the
# sheet's scroll/selection state to reflect the newly added rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source workbook was put into manual calculation mode as part of this edit.
$excel.Calculation = -4135   # xlCalculationManual

# New rows to append, mirroring the existing "superadmin" rows above them
# (regcntr_id, machine_id, device_id, lang_code, is_active, cr_by, cr_dtimes, eff_dtimes)
$newRows = @(
    @(10002, 10032, 3000176, "eng", $true, "superadmin", "now()", "now()"),
    @(10002, 10032, 3000177, "eng", $true, "superadmin", "now()", "now()"),
    @(10002, 10032, 3000178, "eng", $true, "superadmin", "now()", "now()"),
    @(10002, 10032, 3000179, "eng", $true, "superadmin", "now()", "now()"),
    @(10002, 10032, 3000180, "eng", $true, "superadmin", "now()", "now()")
)

$startRow = 157
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $newRows[$i]
    $row = $startRow + $i
    for ($c = 0; $c -lt $r.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $r[$c]
    }
}

# Update the view so the newly entered row/cell is visible and selected,
# matching the author's final cursor position after data entry.
$ws.Activate() | Out-Null
$ws.Range("A151").Select() | Out-Null
$ws.Range("D157").Select() | Out-Null
